$d = $word.ActiveDocument
$replacements = @(
    ,@("2023-03-19 Sunday", "2023-03-20 Monday")
    ,@("81×96=", "43×61=")
    ,@("36×55=", "75×58=")
    ,@("24×89=", "43×43=")
    ,@("39×46=", "35×55=")
    ,@("43×33=", "36×81=")
    ,@("27×21=", "94×19=")
    ,@("21×28=", "15×29=")
    ,@("78×94=", "61×73=")
    ,@("57×29=", "60×91=")
    ,@("53×54=", "47×88=")
    ,@("66×86=", "29×37=")
    ,@("56×52=", "22×41=")
    ,@("19×28=", "79×66=")
    ,@("84×32=", "50×47=")
    ,@("22×47=", "11×58=")
    ,@("94×53=", "77×88=")
    ,@("45×28=", "66×90=")
    ,@("23×36=", "60×21=")
    ,@("41×98=", "54×56=")
    ,@("81×19=", "90×73=")
    ,@("48×30=", "35×10=")
    ,@("26×90=", "38×57=")
    ,@("12×30=", "31×39=")
    ,@("65×89=", "49×43=")
    ,@("90×42=", "22×63=")
    ,@("23×47=", "100×27=")
    ,@("21×96=", "67×91=")
    ,@("84×85=", "78×50=")
    ,@("65×91=", "11×40=")
    ,@("30×66=", "81×98=")
    ,@("19×50=", "48×58=")
    ,@("67×26=", "19×99=")
    ,@("75×17=", "54×88=")
    ,@("20×57=", "14×24=")
    ,@("36×52=", "77×39=")
    ,@("42×64=", "10×29=")
    ,@("100×98=", "45×43=")
    ,@("59×42=", "77×64=")
    ,@("92×94=", "63×85=")
    ,@("87×66=", "28×28=")
    ,@("33×67=", "54×60=")
    ,@("30×94=", "36×39=")
    ,@("53×61=", "36×73=")
    ,@("64×77=", "72×85=")
    ,@("40×60=", "92×46=")
    ,@("44×10=", "15×65=")
    ,@("58×54=", "31×49=")
    ,@("58×58=", "61×92=")
    ,@("38×14=", "66×40=")
    ,@("32×56=", "63×62=")
    ,@("84×45=", "37×85=")
    ,@("59×23=", "45×47=")
    ,@("81×41=", "90×36=")
    ,@("53×12=", "72×45=")
    ,@("47×94=", "85×56=")
    ,@("70×45=", "26×33=")
    ,@("13×64=", "52×70=")
    ,@("23×54=", "77×92=")
    ,@("30×21=", "33×95=")
    ,@("60×11=", "34×66=")
    ,@("97×97=", "68×71=")
    ,@("78×60=", "98×53=")
    ,@("73×83=", "65×50=")
    ,@("40×53=", "36×62=")
    ,@("17×73=", "64×76=")
    ,@("90×35=", "63×30=")
    ,@("11×94=", "51×75=")
    ,@("67×39=", "16×99=")
    ,@("87×70=", "76×59=")
    ,@("96×36=", "24×63=")
    ,@("71×70=", "18×88=")
    ,@("96×86=", "57×59=")
    ,@("29×73=", "86×93=")
    ,@("43×30=", "96×43=")
    ,@("30×14=", "21×11=")
    ,@("14×38=", "36×80=")
    ,@("37×47=", "12×26=")
    ,@("55×67=", "44×75=")
    ,@("20×27=", "17×45=")
    ,@("18×13=", "29×58=")
    ,@("24×79=", "84×79=")
    ,@("11×51=", "17×36=")
    ,@("17×72=", "37×13=")
    ,@("18×22=", "78×30=")
    ,@("33×89=", "17×25=")
    ,@("23×87=", "71×27=")
    ,@("15×87=", "98×98=")
    ,@("81×52=", "60×19=")
    ,@("95×35=", "47×26=")
    ,@("90×32=", "21×79=")
    ,@("74×89=", "59×60=")
    ,@("10×58=", "83×51=")
    ,@("38×59=", "65×55=")
    ,@("63×31=", "37×22=")
    ,@("78×84=", "59×26=")
    ,@("66×95=", "54×59=")
    ,@("61×42=", "74×75=")
    ,@("93×39=", "95×10=")
    ,@("68×65=", "63×45=")
    ,@("34×98=", "29×20=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $rng = $d.Content
    $found = $rng.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Output "NOT FOUND: $old"
    }
}
Write-Output "Replacements complete"
